$wb = $excel.ActiveWorkbook

# Fix typo in worksheet name: CONDUCTOR_opertation -> CONDUCTOR_operation
$ws = $wb.Worksheets.Item("CONDUCTOR_opertation")
$ws.Name = "CONDUCTOR_operation"

# Select cell G5 on the (now renamed) CONDUCTOR_operation sheet and activate that sheet/tab
$ws.Activate()
$ws.Range("G5").Select()
